$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 24
$ws.Range("D2").Value = 15

$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 21

$ws.Range("B4").Value = 26
$ws.Range("C4").Value = 15
$ws.Range("D4").Value = 12

$ws.Range("B5").Value = "Yu Qiao"
$ws.Range("C5").Value = "Luc Van Gool"
$ws.Range("D5").Value = "Lei Zhang"
